$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after "line6" (row 7), pushing the "extr*" rows
# down by two rows. Excel's row Insert() repeats on the same row index so
# that two blank rows end up at rows 8 and 9.
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

# Populate the new row 8 -> line7
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# Populate the new row 9 -> line8
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Match the formatting used by the rest of column A (bold, thin border,
# centered/top aligned) for the two freshly-inserted rows.
$newIndexCells = $ws.Range("A8:A9")
$newIndexCells.Font.Bold = $true
$newIndexCells.HorizontalAlignment = -4108
$newIndexCells.VerticalAlignment = -4160
$newIndexCells.Borders.LineStyle = 1

# The "extr5" row (now row 14, after the two-row shift) flips its
# "in_service" flag from FALSE to TRUE.
$ws.Cells.Item(14, 5).Value = $true

# The index column (A) for every "extr*" row (now rows 10-17) is
# renumbered to keep counting on from the newly inserted rows.
for ($r = 10; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
